$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (shared strings): *_old -> *_FV2210, *_new -> *_FV2304 ---
$ws.Range("A1").Value2 = "Segmentname_FV2210"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2210"
$ws.Range("C1").Value2 = "Segment_FV2210"
$ws.Range("D1").Value2 = "Datenelement_FV2210"
$ws.Range("E1").Value2 = "Segment ID_FV2210"
$ws.Range("F1").Value2 = "Code_FV2210"
$ws.Range("G1").Value2 = "Qualifier_FV2210"
$ws.Range("H1").Value2 = "Beschreibung_FV2210"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value2 = "Bedingung_FV2210"

$ws.Range("L1").Value2 = "Segmentname_FV2304"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2304"
$ws.Range("N1").Value2 = "Segment_FV2304"
$ws.Range("O1").Value2 = "Datenelement_FV2304"
$ws.Range("P1").Value2 = "Segment ID_FV2304"
$ws.Range("Q1").Value2 = "Code_FV2304"
$ws.Range("R1").Value2 = "Qualifier_FV2304"
$ws.Range("S1").Value2 = "Beschreibung_FV2304"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value2 = "Bedingung_FV2304"

# --- Turn the data range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U57"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
